$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Renumber existing rows 9-14 (A column) - shift numbering from 5..10 down to 1..6
$ws.Range("A9").Value = 1
$ws.Range("A10").Value = 2
$ws.Range("A11").Value = 3
$ws.Range("A12").Value = 4
$ws.Range("A13").Value = 5
$ws.Range("A14").Value = 6

# New row 15: project "jbal"
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "jbal"
$ws.Range("C15").Value = 102
$ws.Range("D15").Value = 109
$ws.Range("E15").Value = 113
$ws.Range("F15").Value = 370
$ws.Range("G15").Value = 28
$ws.Range("H15").Value = 250
$ws.Range("I15").Value = 4678
$ws.Range("J15").Value = 4700
$ws.Range("K15").Value = 88
$ws.Range("L15").Value = "java"

# Hyperlink for M15 (display text is the URL itself, matching the other rows)
$ws.Hyperlinks.Add($ws.Cells.Item(15, 13), "https://github.com/datazuul/jbal")

# Match the style used by the other hyperlink cells in column M (M9:M14)
$ws.Range("M14").Copy()
$ws.Range("M15").PasteSpecial(-4122)

# Old row 15 (now row 16) — just the A value
$ws.Range("A16").Value = 8

# Update selection
$ws.Range("E21").Select()
